$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers (bold, matches existing A1/B1 style) ---
$ws.Range("A1").Value = "file_name"
$ws.Range("B1").Value = "experiment_date"
$ws.Range("C1").Value = "description"
$ws.Range("C1").Font.Bold = $true

# --- Row 2 ---
$ws.Range("A2").Value = "S_aureus_30jan2020_N.xlsx"
$ws.Range("B2").Value = "'30-1-2020"
$ws.Range("C2").Value = "S.aureus grown with red and white rose extracts"

# --- Row 3 ---
$ws.Range("A3").Value = "K_pneumoniae_2maart2020_N.xls"
$ws.Range("B3").Value = "'2-3-2020"
$ws.Range("C3").Value = "K.pneumoniae grown with red and white rose extracts"

# --- Row 4 ---
$ws.Range("A4").Value = "S_aureus_stampersruw_12nov2021_N.xls"
$ws.Range("B4").Value = "'12-11-2021"
$ws.Range("C4").Value = "S.aureus grown with tulip stamen extracts"

# --- Row 5 ---
$ws.Range("A5").Value = "E_coli_stampersruw_02dec2021_N.xlsx"
$ws.Range("B5").Value = "'2-12-2021"
$ws.Range("C5").Value = "E.coli grown with red tulip stamen extracts"

# --- Column widths (closest achievable values to the target 38.6640625 /
#     18.6640625 / 48 "characters" widths; the engine quantizes ColumnWidth
#     to 1/6-character steps, so these inputs round-trip to the nearest
#     representable width) ---
$ws.Columns.Item(1).ColumnWidth = 37.833333333333336
$ws.Columns.Item(2).ColumnWidth = 17.833333333333332
$ws.Columns.Item(3).ColumnWidth = 47.166666666666664

# --- Selection ---
$ws.Range("B9").Select() | Out-Null
